$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.968.83"
$ws.Range("E2").Value = "  -3.61%  "
$ws.Range("D3").Value = "1.642.30"
$ws.Range("E3").Value = "  -5.69%  "
$ws.Range("D4").Value = "'0.9987"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'236.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.03%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "'0.4803"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -6.23%  "
$ws.Range("D8").Value = "'0.2587"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.79%  "
$ws.Range("D9").Value = "'0.05999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.01%  "
$ws.Range("D10").Value = "'0.07200"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("D11").Value = "1.644.31"
$ws.Range("E11").Value = "  -5.52%  "
$ws.Range("D12").Value = "'14.84"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.87%  "
$ws.Range("D13").Value = "'0.6192"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.51%  "
$ws.Range("D14").Value = "'4.525"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.24%  "
$ws.Range("D15").Value = "'72.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.28%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "'0.9985"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").Value = "24.977.08"
$ws.Range("E18").Value = "  -3.66%  "
$ws.Range("D19").Value = "'11.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.93%  "
$ws.Range("D20").Value = "'0.000006611"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.98%  "
$ws.Range("D21").Value = "'4.501"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.16%  "
$ws.Range("D22").Value = "1.856.01"
$ws.Range("E22").Value = "  -5.55%  "
$ws.Range("D23").Value = "'8.592"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.68%  "
$ws.Range("D24").Value = "'5.291"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.01%  "
$ws.Range("D25").Value = "'131.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.98%  "
$ws.Range("D26").Value = "'14.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.94%  "
$ws.Range("D27").Value = "'1.394"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.49%  "
$ws.Range("D28").Value = "'102.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.81%  "
$ws.Range("D29").Value = "'1.659"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.47%  "
$ws.Range("D30").Value = "'3.737"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.49%  "
$ws.Range("D31").Value = "'0.07814"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.05%  "
$ws.Range("D32").Value = "'3.549"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.72%  "
$ws.Range("D33").Value = "'0.04426"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.63%  "
$ws.Range("D34").Value = "'0.9998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "'2.585"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.76%  "
$ws.Range("D36").Value = "'0.9283"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.03%  "
$ws.Range("D37").Value = "'0.5854"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.47%  "
$ws.Range("D38").Value = "'2.572"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.69%  "
$ws.Range("E39").Value = "  -2.11%  "
$ws.Range("D40").Value = "'0.8466"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.77%  "
$ws.Range("D41").Value = "'0.9985"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").Value = "'1.808"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.65%  "
$ws.Range("D43").Value = "'97.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.84%  "
$ws.Range("D44").Value = "'0.3716"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.35%  "
$ws.Range("D45").Value = "'4.768"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.55%  "
$ws.Range("D46").Value = "'0.1151"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.86%  "
$ws.Range("D47").Value = "'6.097"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.19%  "
$ws.Range("D48").Value = "'0.05195"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("D49").Value = "'29.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.07%  "
$ws.Range("D50").Value = "'0.9991"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.39%  "

# Row 51: Aave -> USDD (full row replaced)
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'1.000"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.49%  "
